# Apply "changed the SAM site repository structure" edit:
# Add a new data row (row 10) to sheet1 with a new investigator record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new row 10 with values (matching existing row layout / headers)
$ws.Range("A10").Value = "Principal"
$ws.Range("B10").Value = "0000/0011"
$ws.Range("D10").Value = "Robert Jones MD"
$ws.Range("G10").Value = "Robert"
$ws.Range("I10").Value = "Jones"
$ws.Range("J10").Value = "St John's"
$ws.Range("K10").Value = "Bangalore"
$ws.Range("P10").Value = "India"

# Match style used by other cells in column I (center/middle aligned, same as I2:I9)
$ws.Range("I10").HorizontalAlignment = -4108  # xlCenter
$ws.Range("I10").VerticalAlignment = -4108    # xlCenter

# Update selection to reflect the new active cell, as in the diff
$ws.Range("A10").Select()
